$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.701.89'
$ws.Range('E2').Value = '  -1.28%  '
$ws.Range('D3').Value = '2.903.27'
$ws.Range('E3').Value = '  -1.88%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '528.25'
$ws.Range('E5').Value = '  -2.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.79'
$ws.Range('E6').Value = '  -4.98%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.555'
$ws.Range('E8').Value = '  -2.58%  '
$ws.Range('D9').Value = '2.911.49'
$ws.Range('E9').Value = '  -1.83%  '
$ws.Range('E10').Value = '  -4.22%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.04'
$ws.Range('E11').Value = '  -1.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.361'
$ws.Range('E12').Value = '  -1.80%  '
$ws.Range('D13').Value = '3.410.69'
$ws.Range('E13').Value = '  -1.89%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.126'
$ws.Range('E14').Value = '  +1.77%  '
$ws.Range('D15').Value = '60.658.12'
$ws.Range('E15').Value = '  -1.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.84'
$ws.Range('E16').Value = '  -3.42%  '
$ws.Range('D17').Value = '2.914.47'
$ws.Range('E17').Value = '  -1.73%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0000141'
$ws.Range('E18').Value = '  -3.60%  '
$ws.Range('E19').Value = '  -2.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.73'
$ws.Range('E20').Value = '  -2.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '361.79'
$ws.Range('E21').Value = '  -5.33%  '
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.70'
$ws.Range('E24').Value = '  +0.59%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '64.73'
$ws.Range('E25').Value = '  -1.02%  '
$ws.Range('E26').Value = '  -3.01%  '
$ws.Range('E27').Value = '  -3.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.997'
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D30').Value = '0.0₃0854'
$ws.Range('E30').Value = '  -8.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('E32').Value = '  -2.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.79'
$ws.Range('E33').Value = '  -3.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '151.86'
$ws.Range('E34').Value = '  -4.80%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.36'
$ws.Range('E35').Value = '  -5.99%  '
$ws.Range('E36').Value = '  -5.87%  '
$ws.Range('E37').Value = '  -5.64%  '
$ws.Range('E38').Value = '  -4.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.84'
$ws.Range('E39').Value = '  +1.70%  '
$ws.Range('E40').Value = '  -4.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.72'
$ws.Range('E41').Value = '  -5.29%  '
$ws.Range('D42').Value = '2.297.26'
$ws.Range('E42').Value = '  -4.57%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.649'
$ws.Range('E43').Value = '  -2.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0586'
$ws.Range('E44').Value = '  -0.99%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '20.46'
$ws.Range('E45').Value = '  -7.46%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.997'
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.99'
$ws.Range('E47').Value = '  +0.51%  '
$ws.Range('E48').Value = '  -3.13%  '
$ws.Range('B49').Value = 'WhiteBITCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.31'
$ws.Range('E49').Value = '  -1.33%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0925'
$ws.Range('E50').Value = '  -3.45%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '251.16'
$ws.Range('E51').Value = '  -6.21%  '
